$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# The sheet gets two new rows of "text bag" results inserted right after the
# existing "addFeature" block (rows 12-15), pushing the final "emotion.csv"
# summary row from row 17 down to row 19 (row 18 stays blank, same as the
# old gap before row 17).
#
# Step 1: relocate the existing row 17 ("emotion.csv") content down to row 19
# *before* writing anything new into row 16/17, so the shared-string entries
# referenced only by that row move with it instead of being overwritten.
# ---------------------------------------------------------------------------

$ws.Range("A19").Value = "emotion.csv"
$ws.Range("B19").Value = 0.92717749999999999
$ws.Range("C19").Value = 0.92286999999999997
$ws.Range("D19").Value = 0.89654999999999996
$ws.Range("F19").Value = "整段文字的emotion"

$ws.Range("B19").NumberFormat = "0.00000_ "

# Clear out the old row 17 cells now that their content lives on row 19.
$ws.Range("A17").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("F17").ClearContents()

# ---------------------------------------------------------------------------
# Step 2: write the two new "text bag" score rows (16 and 17).
# ---------------------------------------------------------------------------

$ws.Range("B16").Value = 0.94431229999999999
$ws.Range("C16").Value = 0.91515000000000002
$ws.Range("D16").Value = 0.89575000000000005
$ws.Range("F16").Value = "text bag的方式，增加属性"

$ws.Range("B17").Value = 0.92886159999999995
$ws.Range("C17").Value = 0.92737000000000003
$ws.Range("D17").Value = 0.89566000000000001
$ws.Range("F17").Value = "headline, snippet, abstract都通过text bag方式，增加属性"

$ws.Range("B16:D17").NumberFormat = "0.00000_ "

# ---------------------------------------------------------------------------
# Step 3: cosmetics matching the author's final selection state.
# ---------------------------------------------------------------------------

[void]$ws.Range("D17").Select()
